$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:M25 (columns B..M, rows 2..25)
$data = New-Object 'object[,]' 24,12

$data[0,0] = 0.6242483216389303
$data[0,1] = 0.04711193483127829
$data[0,2] = 0.6218071849043838
$data[0,3] = 0.2328265924522839
$data[0,4] = 0
$data[0,5] = 3.343399241024144
$data[0,6] = 2.492473251076007
$data[0,7] = 0
$data[0,8] = 0.1057256128926127
$data[0,9] = 0.5928307476229691
$data[0,10] = 0
$data[0,11] = 0.4046693509676231
$data[1,0] = 0.6038341079080283
$data[1,1] = 0.04536461164018135
$data[1,2] = 0.6164722342893185
$data[1,3] = 0.2311791786835542
$data[1,4] = 0
$data[1,5] = 3.270865788480393
$data[1,6] = 2.461904748980345
$data[1,7] = 0
$data[1,8] = 0.105205887210257
$data[1,9] = 0.5731166599070718
$data[1,10] = 0
$data[1,11] = 0.3979037422451697
$data[2,0] = 0.591822810942773
$data[2,1] = 0.04435308912982805
$data[2,2] = 0.6134784389149104
$data[2,3] = 0.2302649964697849
$data[2,4] = 0
$data[2,5] = 3.227099439017593
$data[2,6] = 2.443663051185041
$data[2,7] = 0
$data[2,8] = 0.1049244909895748
$data[2,9] = 0.561552292982384
$data[2,10] = 0
$data[2,11] = 0.3940231699484755
$data[3,0] = 0.5870595386332411
$data[3,1] = 0.04395625349270915
$data[3,2] = 0.6123293235247331
$data[3,3] = 0.2299169340166536
$data[3,4] = 0
$data[3,5] = 3.209456750158637
$data[3,6] = 2.436361672836341
$data[3,7] = 0
$data[3,8] = 0.1048193024060069
$data[3,9] = 0.5569753729473916
$data[3,10] = 0
$data[3,11] = 0.3925105569239804
$data[4,0] = 0.5862765357636306
$data[4,1] = 0.04389128579252599
$data[4,2] = 0.6121427948087046
$data[4,3] = 0.2298606166624566
$data[4,4] = 0
$data[4,5] = 3.206538788667586
$data[4,6] = 2.435157261816585
$data[4,7] = 0
$data[4,8] = 0.1048024086091459
$data[4,9] = 0.5562235670655298
$data[4,10] = 0
$data[4,11] = 0.3922635402619505
$data[5,0] = 0.5917580397816948
$data[5,1] = 0.04434767511605742
$data[5,2] = 0.6134626545510429
$data[5,3] = 0.2302602032840717
$data[5,4] = 0
$data[5,5] = 3.226860725689761
$data[5,6] = 2.443564047243058
$data[5,7] = 0
$data[5,8] = 0.1049230339867329
$data[5,9] = 0.5614900180109714
$data[5,10] = 0
$data[5,11] = 0.3940024919875853
$data[6,0] = 0.6171008502637108
$data[6,1] = 0.04649669124367506
$data[6,2] = 0.6199091835094777
$data[6,3] = 0.2322383596436879
$data[6,4] = 0
$data[6,5] = 3.318229348915395
$data[6,6] = 2.481823511231596
$data[6,7] = 0
$data[6,8] = 0.1055385817606265
$data[6,9] = 0.5859211033101133
$data[6,10] = 0
$data[6,11] = 0.402279772569436
$data[7,0] = 0.6709578984147697
$data[7,1] = 0.05120056582471477
$data[7,2] = 0.6347890303572399
$data[7,3] = 0.2368904115495383
$data[7,4] = 0
$data[7,5] = 3.503565833452001
$data[7,6] = 2.561057706534427
$data[7,7] = 0
$data[7,8] = 0.1070451891482449
$data[7,9] = 0.6381283138268259
$data[7,10] = 0
$data[7,11] = 0.4206850046007702
$data[8,0] = 0.7130813828496514
$data[8,1] = 0.05495965706795403
$data[8,2] = 0.6470900958042876
$data[8,3] = 0.2407809504866023
$data[8,4] = 0
$data[8,5] = 3.643584228727832
$data[8,6] = 2.621873998738494
$data[8,7] = 0
$data[8,8] = 0.1083352835350979
$data[8,9] = 0.6791275767423315
$data[8,10] = 0
$data[8,11] = 0.435539062135021
$data[9,0] = 0.7328037067759112
$data[9,1] = 0.05673664772950815
$data[9,2] = 0.6529845206261768
$data[9,3] = 0.2426538771520512
$data[9,4] = 0
$data[9,5] = 3.70813995142845
$data[9,6] = 2.650114773607299
$data[9,7] = 0
$data[9,8] = 0.1089621079877077
$data[9,9] = 0.6983583688295596
$data[9,10] = 0
$data[9,11] = 0.4425873798903197
$data[10,0] = 0.7403528468807679
$data[10,1] = 0.05741925415473759
$data[10,2] = 0.6552595799144285
$data[10,3] = 0.2433779507592035
$data[10,4] = 0
$data[10,5] = 3.732710695750143
$data[10,6] = 2.660892079102894
$data[10,7] = 0
$data[10,8] = 0.1092052234326033
$data[10,9] = 0.705724308003056
$data[10,10] = 0
$data[10,11] = 0.4452983523398473
$data[11,0] = 0.7387234113617467
$data[11,1] = 0.05727181017375926
$data[11,2] = 0.6547676936616256
$data[11,3] = 0.2432213484891506
$data[11,4] = 0
$data[11,5] = 3.727413373545176
$data[11,6] = 2.658567289045664
$data[11,7] = 0
$data[11,8] = 0.1091526083230576
$data[11,9] = 0.7041341965384902
$data[11,10] = 0
$data[11,11] = 0.4447126301400104
$data[12,0] = 0.7334231602038983
$data[12,1] = 0.05679261139172809
$data[12,2] = 0.6531708297936518
$data[12,3] = 0.2427131497022188
$data[12,4] = 0
$data[12,5] = 3.710158891788296
$data[12,6] = 2.650999760333264
$data[12,7] = 0
$data[12,8] = 0.1089819939589205
$data[12,9] = 0.6989626910370248
$data[12,10] = 0
$data[12,11] = 0.4428095725624885
$data[13,0] = 0.7301871213678339
$data[13,1] = 0.05650035344200433
$data[13,2] = 0.6521983012464432
$data[13,3] = 0.2424037958236838
$data[13,4] = 0
$data[13,5] = 3.699606322037482
$data[13,6] = 2.646375270299245
$data[13,7] = 0
$data[13,8] = 0.1088782367989296
$data[13,9] = 0.6958058965683449
$data[13,10] = 0
$data[13,11] = 0.441649356973393
$data[14,0] = 0.7118037689550931
$data[14,1] = 0.05484487990379705
$data[14,2] = 0.6467108924599358
$data[14,3] = 0.2406606259098396
$data[14,4] = 0
$data[14,5] = 3.639382775555049
$data[14,6] = 2.620040007821075
$data[14,7] = 0
$data[14,8] = 0.1082951235045186
$data[14,9] = 0.6778824918633291
$data[14,10] = 0
$data[14,11] = 0.4350843016263752
$data[15,0] = 0.7006697659845997
$data[15,1] = 0.05384650000532076
$data[15,2] = 0.6434210454270612
$data[15,3] = 0.239617661887138
$data[15,4] = 0
$data[15,5] = 3.602658817108619
$data[15,6] = 2.604031790406452
$data[15,7] = 0
$data[15,8] = 0.107947639229053
$data[15,9] = 0.6670357954889994
$data[15,10] = 0
$data[15,11] = 0.4311314625111464
$data[16,0] = 0.6943184938877494
$data[16,1] = 0.05327855674399018
$data[16,2] = 0.6415569169627702
$data[16,3] = 0.2390274801982102
$data[16,4] = 0
$data[16,5] = 3.581617107984187
$data[16,6] = 2.594878395758172
$data[16,7] = 0
$data[16,8] = 0.1077515357922749
$data[16,9] = 0.660851641561635
$data[16,10] = 0
$data[16,11] = 0.4288852927865889
$data[17,0] = 0.6921771103567096
$data[17,1] = 0.05308734057854281
$data[17,2] = 0.6409305812368871
$data[17,3] = 0.2388293211028056
$data[17,4] = 0
$data[17,5] = 3.574506614056389
$data[17,6] = 2.591788493695702
$data[17,7] = 0
$data[17,8] = 0.1076857842222125
$data[17,9] = 0.6587671598393854
$data[17,10] = 0
$data[17,11] = 0.4281294822172939
$data[18,0] = 0.7018495428768006
$data[18,1] = 0.05395212688469542
$data[18,2] = 0.6437683463454391
$data[18,3] = 0.2397276828213961
$data[18,4] = 0
$data[18,5] = 3.606559759989779
$data[18,6] = 2.605730289361901
$data[18,7] = 0
$data[18,8] = 0.1079842403062941
$data[18,9] = 0.6681847943530954
$data[18,10] = 0
$data[18,11] = 0.4315494128706661
$data[19,0] = 0.7349777801722155
$data[19,1] = 0.05693309990984119
$data[19,2] = 0.6536387010213218
$data[19,3] = 0.2428620172414711
$data[19,4] = 0
$data[19,5] = 3.715223552700479
$data[19,6] = 2.653220267168081
$data[19,7] = 0
$data[19,8] = 0.1090319514192331
$data[19,9] = 0.7004794151423539
$data[19,10] = 0
$data[19,11] = 0.4433674082841534
$data[20,0] = 0.7570996078193843
$data[20,1] = 0.05893790252024189
$data[20,2] = 0.6603400173897853
$data[20,3] = 0.244996972781685
$data[20,4] = 0
$data[20,5] = 3.786970127882398
$data[20,6] = 2.684742503163022
$data[20,7] = 0
$data[20,8] = 0.1097502129568682
$data[20,9] = 0.7220735475856941
$data[20,10] = 0
$data[20,11] = 0.4513355687225413
$data[21,0] = 0.7452496507533795
$data[21,1] = 0.05786270258209925
$data[21,2] = 0.6567404704615569
$data[21,3] = 0.2438495893558894
$data[21,4] = 0
$data[21,5] = 3.748610595684625
$data[21,6] = 2.667873989990028
$data[21,7] = 0
$data[21,8] = 0.1093637941780727
$data[21,9] = 0.710503640584335
$data[21,10] = 0
$data[21,11] = 0.4470604255925821
$data[22,0] = 0.7013160105924499
$data[22,1] = 0.05390435416794048
$data[22,2] = 0.6436112467290798
$data[22,3] = 0.2396779129785642
$data[22,4] = 0
$data[22,5] = 3.604795921445344
$data[22,6] = 2.604962242436159
$data[22,7] = 0
$data[22,8] = 0.1079676815287272
$data[22,9] = 0.6676651707724375
$data[22,10] = 0
$data[22,11] = 0.4313603753794197
$data[23,0] = 0.6559408890243503
$data[23,1] = 0.04987514661300452
$data[23,2] = 0.6305235796982345
$data[23,3] = 0.2355490130001634
$data[23,4] = 0
$data[23,5] = 3.452757873053571
$data[23,6] = 2.539168411983695
$data[23,7] = 0
$data[23,8] = 0.1066054859843888
$data[23,9] = 0.6235424649302672
$data[23,10] = 0
$data[23,11] = 0.4154725323017274

$ws.Range("B2:M25").Value = $data